$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 30-37 (table shrinks from 35 to 27 data rows)
$ws.Range("A30:F37").EntireRow.Delete() | Out-Null

# Update changed cell contents for rows 2-29 (columns B-F; column A index numbers unchanged)
$ws.Range("B2").Value2 = "NSE:ADROITINFO"
$ws.Range("C2").Value2 = "NSE:ALKALI"
$ws.Range("D2").Value2 = "NSE:BAJAJFINSV"
$ws.Range("E2").Value2 = "NSE:ADANIENT"
$ws.Range("F2").Value2 = "NSE:BOSCHLTD"
$ws.Range("B3").Value2 = "NSE:AGSTRA"
$ws.Range("C3").Value2 = "NSE:CANTABIL"
$ws.Range("E3").Value2 = ""
$ws.Range("F3").Value2 = "NSE:DIXON"
$ws.Range("B4").Value2 = "NSE:CENTURYTEX"
$ws.Range("C4").Value2 = "NSE:CHAMBLFERT"
$ws.Range("E4").Value2 = ""
$ws.Range("F4").Value2 = "NSE:JSWSTEEL"
$ws.Range("B5").Value2 = "NSE:CONFIPET"
$ws.Range("C5").Value2 = "NSE:IEX"
$ws.Range("E5").Value2 = ""
$ws.Range("B6").Value2 = "NSE:CRAFTSMAN"
$ws.Range("C6").Value2 = "NSE:MEGASOFT"
$ws.Range("E6").Value2 = ""
$ws.Range("B7").Value2 = "NSE:GILLETTE"
$ws.Range("C7").Value2 = "NSE:NDL"
$ws.Range("E7").Value2 = ""
$ws.Range("B8").Value2 = "NSE:HDFCNEXT50"
$ws.Range("C8").Value2 = "NSE:NEOGEN"
$ws.Range("E8").Value2 = ""
$ws.Range("B9").Value2 = "NSE:HILTON"
$ws.Range("C9").Value2 = "NSE:PGHL"
$ws.Range("B10").Value2 = "NSE:HNGSNGBEES"
$ws.Range("C10").Value2 = "NSE:RANASUG"
$ws.Range("B11").Value2 = "NSE:HONDAPOWER"
$ws.Range("B12").Value2 = "NSE:INFRABEES"
$ws.Range("B13").Value2 = "NSE:JETFREIGHT"
$ws.Range("B14").Value2 = "NSE:JHS"
$ws.Range("B15").Value2 = "NSE:JISLDVREQS"
$ws.Range("B16").Value2 = "NSE:JSWSTEEL"
$ws.Range("B17").Value2 = "NSE:JYOTISTRUC"
$ws.Range("B18").Value2 = "NSE:KOKUYOCMLN"
$ws.Range("B19").Value2 = "NSE:MAHSEAMLES"
$ws.Range("B20").Value2 = "NSE:MANGALAM"
$ws.Range("B21").Value2 = "NSE:MOLOWVOL"
$ws.Range("B22").Value2 = "NSE:MVGJL"
$ws.Range("B23").Value2 = "NSE:NATIONALUM"
$ws.Range("B24").Value2 = "NSE:NV20BEES"
$ws.Range("B25").Value2 = "NSE:PARACABLES"
$ws.Range("B26").Value2 = "NSE:PRAKASH"
$ws.Range("B27").Value2 = "NSE:RAIN"
$ws.Range("B28").Value2 = "NSE:ROUTE"
$ws.Range("B29").Value2 = "NSE:SAMBHAAV"
